# The commit adds one new price-observation row for "Haba" (Mercado Mayorista
# Lo Valledor de Santiago) that belongs chronologically right before the
# existing row 158, pushing the former rows 158..271 down to 159..272 and
# growing the sheet's used range from A1:R271 to A1:R272.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 158; Excel automatically shifts every
# row from 158 downward (previously 158..271) to 159..272 and extends the
# sheet dimension accordingly.
$ws.Rows("158:158").Insert()

# Populate the newly inserted row 158 with the new record's data.
$ws.Range("A158").Value = 6
$ws.Range("B158").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C158").Value = "Metropolitana"
$ws.Range("D158").Value = 44767
$ws.Range("E158").Value = 13
$ws.Range("F158").Value = 100112026
$ws.Range("G158").Value = "Haba"
$ws.Range("H158").Value = "Sin especificar"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 550
$ws.Range("K158").Value = 13000
$ws.Range("L158").Value = 15000
$ws.Range("M158").Value = 14091
$ws.Range("N158").Value = "`$/saco 25 kilos"
$ws.Range("O158").Value = "Región de Coquimbo"
$ws.Range("P158").Value = 564
$ws.Range("Q158").Value = 25
$ws.Range("R158").Value = "Hortaliza"
